$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 11 de Mayo de 2020 a las 23:35"

# Update changed numeric values in the data table (rows 4-69)
$ws.Range("B4").Value = 64853
$ws.Range("C4").Value = 39604
$ws.Range("D4").Value = 16566
$ws.Range("E4").Value = 8683
$ws.Range("B5").Value = 54807
$ws.Range("C5").Value = 24454
$ws.Range("D5").Value = 24798
$ws.Range("E5").Value = 5555
$ws.Range("B6").Value = 17995
$ws.Range("C6").Value = 7453
$ws.Range("D6").Value = 8637
$ws.Range("E6").Value = 1905
$ws.Range("B7").Value = 16387
$ws.Range("C7").Value = 6172
$ws.Range("D7").Value = 7429
$ws.Range("E7").Value = 2786
$ws.Range("B9").Value = 12341
$ws.Range("C9").Value = 9205
$ws.Range("D9").Value = 1814
$ws.Range("E9").Value = 1322
$ws.Range("B10").Value = 10021
$ws.Range("E10").Value = 790
$ws.Range("B11").Value = 9280
$ws.Range("C11").Value = 7445
$ws.Range("D11").Value = 1242
$ws.Range("E11").Value = 593
$ws.Range("B12").Value = 8293
$ws.Range("C12").Value = 1967
$ws.Range("D12").Value = 5300
$ws.Range("E12").Value = 1026
$ws.Range("B14").Value = 5336
$ws.Range("C14").Value = 3116
$ws.Range("D14").Value = 1395
$ws.Range("E14").Value = 825
$ws.Range("B15").Value = 5307
$ws.Range("C15").Value = 1883
$ws.Range("D15").Value = 2708
$ws.Range("E15").Value = 716
$ws.Range("B16").Value = 5065
$ws.Range("C16").Value = 3070
$ws.Range("D16").Value = 1505
$ws.Range("E16").Value = 490
$ws.Range("B17").Value = 4959
$ws.Range("C17").Value = 2347
$ws.Range("D17").Value = 1968
$ws.Range("B18").Value = 4948
$ws.Range("C18").Value = 1335
$ws.Range("D18").Value = 3131
$ws.Range("E18").Value = 482
$ws.Range("B20").Value = 4649
$ws.Range("E20").Value = 351
$ws.Range("B21").Value = 4242
$ws.Range("C21").Value = 1487
$ws.Range("D21").Value = 2408
$ws.Range("E21").Value = 347
$ws.Range("B22").Value = 4012
$ws.Range("C22").Value = 1124
$ws.Range("D22").Value = 2541
$ws.Range("E22").Value = 347
$ws.Range("B23").Value = 4008
$ws.Range("C23").Value = 2700
$ws.Range("D23").Value = 962
$ws.Range("B24").Value = 3760
$ws.Range("C24").Value = 1942
$ws.Range("D24").Value = 1547
$ws.Range("E24").Value = 271
$ws.Range("B25").Value = 3483
$ws.Range("C25").Value = 1535
$ws.Range("D25").Value = 1552
$ws.Range("E25").Value = 396
$ws.Range("B26").Value = 3310
$ws.Range("C26").Value = 837
$ws.Range("D26").Value = 2274
$ws.Range("E26").Value = 199
$ws.Range("B27").Value = 3007
$ws.Range("E27").Value = 277
$ws.Range("B28").Value = 2973
$ws.Range("C28").Value = 2269
$ws.Range("D28").Value = 431
$ws.Range("B29").Value = 2918
$ws.Range("C29").Value = 2409
$ws.Range("D29").Value = 26
$ws.Range("E29").Value = 483
$ws.Range("B30").Value = 2916
$ws.Range("C30").Value = 579
$ws.Range("D30").Value = 2040
$ws.Range("E30").Value = 297
$ws.Range("B31").Value = 2831
$ws.Range("C31").Value = 1449
$ws.Range("D31").Value = 1114
$ws.Range("E31").Value = 268
$ws.Range("B32").Value = 2627
$ws.Range("C32").Value = 867
$ws.Range("D32").Value = 1555
$ws.Range("E32").Value = 205
$ws.Range("B33").Value = 2346
$ws.Range("C33").Value = 1026
$ws.Range("D33").Value = 1021
$ws.Range("E33").Value = 299
$ws.Range("B34").Value = 2262
$ws.Range("C34").Value = 380
$ws.Range("D34").Value = 1764
$ws.Range("E34").Value = 118
$ws.Range("B35").Value = 2260
$ws.Range("C35").Value = 1391
$ws.Range("D35").Value = 720
$ws.Range("E35").Value = 149
$ws.Range("B36").Value = 2245
$ws.Range("C36").Value = 1963
$ws.Range("D36").Value = 81
$ws.Range("E36").Value = 201
$ws.Range("B37").Value = 2198
$ws.Range("C37").Value = 365
$ws.Range("D37").Value = 1595
$ws.Range("E37").Value = 238
$ws.Range("D38").Value = 699
$ws.Range("B41").Value = 1866
$ws.Range("C41").Value = 602
$ws.Range("D41").Value = 1133
$ws.Range("B42").Value = 1654
$ws.Range("C42").Value = 1042
$ws.Range("D42").Value = 447
$ws.Range("E42").Value = 165
$ws.Range("B43").Value = 1589
$ws.Range("C43").Value = 1493
$ws.Range("D43").Value = 0
$ws.Range("B46").Value = 1444
$ws.Range("C46").Value = 469
$ws.Range("D46").Value = 833
$ws.Range("E46").Value = 142
$ws.Range("B47").Value = 1439
$ws.Range("C47").Value = 787
$ws.Range("D47").Value = 550
$ws.Range("E47").Value = 102
$ws.Range("B48").Value = 1155
$ws.Range("C48").Value = 315
$ws.Range("D48").Value = 760
$ws.Range("B49").Value = 1022
$ws.Range("C49").Value = 343
$ws.Range("D49").Value = 581
$ws.Range("E49").Value = 98
$ws.Range("B51").Value = 829
$ws.Range("C51").Value = 306
$ws.Range("D51").Value = 441
$ws.Range("B53").Value = 642
$ws.Range("C53").Value = 403
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 450
$ws.Range("B54").Value = 623
$ws.Range("C54").Value = 347
$ws.Range("D54").Value = 195
$ws.Range("B56").Value = 510
$ws.Range("C56").Value = 303
$ws.Range("D56").Value = 160
$ws.Range("E56").Value = 47
$ws.Range("C59").Value = 113
$ws.Range("D59").Value = 4
$ws.Range("C60").Value = 44
$ws.Range("D60").Value = 45
$ws.Range("E60").Value = 6
$ws.Range("C61").Value = 67
$ws.Range("D61").Value = 10
$ws.Range("E61").Value = 6
$ws.Range("B63").Value = 45
$ws.Range("C63").Value = 40
$ws.Range("D63").Value = 5
$ws.Range("C66").Value = 8
$ws.Range("D66").Value = 0
$ws.Range("C68").Value = 2
$ws.Range("D68").Value = 1
